$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing row 5 raw inputs (C5, D5) ---
# E5/F5 formulas already exist in the workbook and recalculate automatically.
$ws.Range("C5").Value = -10.55
$ws.Range("D5").Value = 11.1

# --- Fill in the raw inputs for the newly populated rows 6-10 ---
$ws.Range("C6").Value = -16.79
$ws.Range("D6").Value = 11.14

$ws.Range("C7").Value = -11.1
$ws.Range("D7").Value = 18.95

$ws.Range("C8").Value = -30.82
$ws.Range("D8").Value = 80.27

$ws.Range("C9").Value = -56.91
$ws.Range("D9").Value = 56.91

$ws.Range("C10").Value = -52.04
$ws.Range("D10").Value = 63.88

# --- Add the offset/gain formulas for the new rows, entered as one fill
# operation per column so Excel groups them into shared formulas, just like
# it would if the user filled the formula down from row 6 through row 10. ---
$ws.Range("E6:E10").Formula = "=`$B`$3-AVERAGE(C6:D6)"
$ws.Range("F6:F10").Formula = "=(`$B`$1-`$B`$2)/(D6-C6)"

# --- Match the print/page setup tweak Excel persisted on save ---
$ws.PageSetup.Orientation = 1

# --- Leave the selection where the user last clicked ---
$ws.Range("F10").Select()
